$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.280.23"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.858.52"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7023"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.16"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07912"
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3032"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.50"
$ws.Range("E10").Value = "  +6.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08165"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.880.06"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.213"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7069"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.51"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.351.08"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.814"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007832"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.30"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.130.10"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.570"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.50"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.900"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.07"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.907"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.400"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.479"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.296"
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.178"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7086"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.004"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.679"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01850"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.688"
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.140.00"
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9211"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.957"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4244"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.38"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.86"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5313"
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.746"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.186"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("E51").Value = "  +0.74%  "
